$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.959.50'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.247.74'
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('E4').Value = '  -0.09%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '232.23'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +0.52%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.641'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +2.21%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '62.98'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -1.78%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +5.72%  '
$ws.Range('E10').Value = '  +3.04%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '57.31'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -0.62%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '26.33'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').Value = '2.581.96'
$ws.Range('E14').Value = '  -1.44%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '15.47'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -1.52%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '6.07'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +3.12%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.828'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +1.90%  '
$ws.Range('D18').Value = '2.257.79'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').Value = '43.850.70'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').Value = '0.0₃0979'
$ws.Range('E20').Value = '  +3.94%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '72.68'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('E22').Value = '  -1.42%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '248.07'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('E24').Value = '  +0.05%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.43'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -5.37%  '
$ws.Range('B26').Value = 'WEMIXToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '3.33'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +21.24%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '2.22'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -4.81%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '9.76'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -0.58%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '173.14'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +1.11%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '21.01'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +2.58%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.139'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +1.25%  '
$ws.Range('E32').Value = '  -1.54%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.123'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +1.02%  '
$ws.Range('E34').Value = '  -1.77%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '4.81'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +2.17%  '
$ws.Range('E36').Value = '  -3.70%  '
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('E39').Value = '  -3.60%  '
$ws.Range('E40').Value = '  +2.20%  '
$ws.Range('E41').Value = '  +0.07%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '8.65'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +1.87%  '
$ws.Range('E43').Value = '  +1.33%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '97.53'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '16.98'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('E47').Value = '  -1.79%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '4.34'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -7.97%  '
$ws.Range('D49').Value = '1.439.20'
$ws.Range('E49').Value = '  -2.81%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '2.30'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -2.05%  '
$ws.Range('E51').Value = '  +1.53%  '
